# Refresh the cryptocurrency price/volume snapshot in the "cryptos" sheet:
# updates the Price (D) and hourly Volume change (E) columns for every
# coin row, and re-seats the Monero / FirstDigitalUSD rows (40 & 41) to
# reflect their new order in the source feed.
#
# Some of the new Price strings (e.g. "610.96") are plain decimals that
# Excel would otherwise auto-coerce to a Number on assignment, unlike the
# thousands-dotted ones (e.g. "67.551.25") which already stay text. Force
# NumberFormat "@" immediately before those writes so the cell keeps the
# same text type as the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.551.25"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "3.522.81"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.96"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.78"
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("D7").Value = "3.521.63"
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.482"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.07"
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.426"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("E13").Value = "  -2.20%  "
$ws.Range("D14").Value = "4.119.59"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.94"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "3.523.81"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").Value = "67.502.31"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.117"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.42"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.25"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "446.26"
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.26"
$ws.Range("E22").Value = "  -4.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.625"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.41"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("E25").Value = "  +10.88%  "
$ws.Range("D26").Value = "3.664.25"
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.20"
$ws.Range("E28").Value = "  -4.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.36"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("E31").Value = "  -3.84%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +4.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.82"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.17"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "3.514.54"
$ws.Range("E36").Value = "  -1.38%  "
$ws.Range("E37").Value = "  -3.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.07"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "177.88"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.21"
$ws.Range("E42").Value = "  +4.21%  "
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("E44").Value = "  -3.20%  "
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.57"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.32"
$ws.Range("E47").Value = "  -4.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.62"
$ws.Range("E48").Value = "  +1.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.26"
$ws.Range("E49").Value = "  +4.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.60"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("E51").Value = "  -0.73%  "
